# New PO forecast model
# Updates the three data sheets (Weekly Quantity, Monthly Trend, PO Forecast)
# with refreshed forecast numbers and appends the newly-forecast weeks/months.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Weekly Quantity" -> append 2 new weekly rows (rows 12-13)
# ---------------------------------------------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")

$weeklyNewRows = @(
    @(45662.99999999999, 11),
    @(45669.99999999999, 6)
)

$startRow = 12
for ($i = 0; $i -lt $weeklyNewRows.Count; $i++) {
    $r = $startRow + $i
    $wsWeekly.Cells.Item($r, 1).Value = $weeklyNewRows[$i][0]
    $wsWeekly.Cells.Item($r, 1).NumberFormat = $wsWeekly.Cells.Item($r - 1, 1).NumberFormat
    $wsWeekly.Cells.Item($r, 2).Value = $weeklyNewRows[$i][1]
}

# ---------------------------------------------------------------------------
# Sheet 2: "Monthly Trend" -> append 1 new monthly row (row 6)
# ---------------------------------------------------------------------------
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

$wsMonthly.Cells.Item(6, 1).Value = 45688.99999999999
$wsMonthly.Cells.Item(6, 1).NumberFormat = $wsMonthly.Cells.Item(5, 1).NumberFormat
$wsMonthly.Cells.Item(6, 2).Value = 17

# ---------------------------------------------------------------------------
# Sheet 3: "PO Forecast" -> refresh forecast values and extend by 2 rows
# ---------------------------------------------------------------------------
$wsForecast = $wb.Worksheets.Item("PO Forecast")

# Final state for rows 2..21 of the forecast sheet (ds, PO_Forecast)
$forecastRows = @(
    @(45466.99999999999, 25),
    @(45480.99999999999, 24),
    @(45487.99999999999, 23),
    @(45494.99999999999, 23),
    @(45501.99999999999, 22),
    @(45508.99999999999, 22),
    @(45515.99999999999, 22),
    @(45522.99999999999, 21),
    @(45634.99999999999, 13),
    @(45641.99999999999, 13),
    @(45662.99999999999, 11),
    @(45669.99999999999, 11),
    @(45676.99999999999, 11),
    @(45683.99999999999, 10),
    @(45690.99999999999, 10),
    @(45697.99999999999, 9),
    @(45704.99999999999, 9),
    @(45711.99999999999, 8),
    @(45718.99999999999, 8),
    @(45725.99999999999, 7)
)

$dateFormat = $wsForecast.Cells.Item(2, 1).NumberFormat

$startRow = 2
for ($i = 0; $i -lt $forecastRows.Count; $i++) {
    $r = $startRow + $i
    $wsForecast.Cells.Item($r, 1).Value = $forecastRows[$i][0]
    $wsForecast.Cells.Item($r, 1).NumberFormat = $dateFormat
    $wsForecast.Cells.Item($r, 2).Value = $forecastRows[$i][1]
}
